$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2311
$ws.Range("E2").Value = 233
$ws.Range("F2").Value = 233
$ws.Range("G2").Value = 257
$ws.Range("H2").Value = 197
$ws.Range("I2").Value = 197
$ws.Range("K2").Value = 3433
$ws.Range("L2").Value = 611
$ws.Range("M2").Value = 2822
$ws.Range("N2").Value = 2822
$ws.Range("P2").Value = 205
$ws.Range("Q2").Value = 188
$ws.Range("R2").Value = -729
$ws.Range("S2").Value = 313
$ws.Range("T2").Value = 135
$ws.Range("U2").Value = 53
$ws.Range("V2").Value = 117
$ws.Range("W2").Value = 10.06
$ws.Range("X2").Value = 8.539999999999999
$ws.Range("Y2").Value = 7.65
$ws.Range("Z2").Value = 6.17
$ws.Range("AA2").Value = 21.63
$ws.Range("AB2").Value = 1277.87
$ws.Range("AC2").Value = 515
$ws.Range("AD2").Value = 11.27
$ws.Range("AE2").Value = 6884
$ws.Range("AF2").Value = 0.84
$ws.Range("AG2").Value = 120
$ws.Range("AH2").Value = 2.07
$ws.Range("AI2").Value = 24.93
$ws.Range("AJ2").Value = 40996887
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 2425
$ws.Range("E3").Value = 255
$ws.Range("F3").Value = 255
$ws.Range("G3").Value = 304
$ws.Range("H3").Value = 227
$ws.Range("I3").Value = 227
$ws.Range("K3").Value = 3474
$ws.Range("L3").Value = 655
$ws.Range("M3").Value = 2819
$ws.Range("N3").Value = 2819
$ws.Range("P3").Value = 205
$ws.Range("Q3").Value = 19
$ws.Range("R3").Value = 432
$ws.Range("S3").Value = -219
$ws.Range("T3").Value = 69
$ws.Range("U3").Value = -50
$ws.Range("V3").Value = 146
$ws.Range("W3").Value = 10.53
$ws.Range("X3").Value = 9.359999999999999
$ws.Range("Y3").Value = 8.039999999999999
$ws.Range("Z3").Value = 6.57
$ws.Range("AA3").Value = 23.25
$ws.Range("AB3").Value = 1361.91
$ws.Range("AC3").Value = 553
$ws.Range("AD3").Value = 11.76
$ws.Range("AE3").Value = 7397
$ws.Range("AF3").Value = 0.88
$ws.Range("AG3").Value = 130
$ws.Range("AH3").Value = 2
$ws.Range("AI3").Value = 21.83
$ws.Range("AJ3").Value = 40996887
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 2576
$ws.Range("E4").Value = 269
$ws.Range("F4").Value = 269
$ws.Range("G4").Value = 245
$ws.Range("H4").Value = 165
$ws.Range("I4").Value = 165
$ws.Range("K4").Value = 3895
$ws.Range("L4").Value = 734
$ws.Range("M4").Value = 3160
$ws.Range("N4").Value = 3160
$ws.Range("P4").Value = 220
$ws.Range("Q4").Value = 600
$ws.Range("R4").Value = -267
$ws.Range("S4").Value = 168
$ws.Range("T4").Value = 125
$ws.Range("U4").Value = 475
$ws.Range("V4").Value = 136
$ws.Range("W4").Value = 10.44
$ws.Range("X4").Value = 6.39
$ws.Range("Y4").Value = 5.51
$ws.Range("Z4").Value = 4.47
$ws.Range("AA4").Value = 23.24
$ws.Range("AB4").Value = 1415.26
$ws.Range("AC4").Value = 379
$ws.Range("AD4").Value = 15.95
$ws.Range("AE4").Value = 7704
$ws.Range("AF4").Value = 0.79
$ws.Range("AG4").Value = 130
$ws.Range("AH4").Value = 2.15
$ws.Range("AI4").Value = 32.38
$ws.Range("AJ4").Value = 43960757
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 2391
$ws.Range("E5").Value = 58
$ws.Range("F5").Value = 58
$ws.Range("G5").Value = 105
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = 7
$ws.Range("K5").Value = 3664
$ws.Range("L5").Value = 561
$ws.Range("M5").Value = 3103
$ws.Range("N5").Value = 3103
$ws.Range("P5").Value = 220
$ws.Range("Q5").Value = -154
$ws.Range("R5").Value = 130
$ws.Range("S5").Value = -53
$ws.Range("T5").Value = 104
$ws.Range("U5").Value = -259
$ws.Range("V5").Value = 136
$ws.Range("W5").Value = 2.44
$ws.Range("X5").Value = 0.31
$ws.Range("Y5").Value = 0.24
$ws.Range("Z5").Value = 0.2
$ws.Range("AA5").Value = 18.08
$ws.Range("AB5").Value = 1395.07
$ws.Range("AC5").Value = 17
$ws.Range("AD5").Value = 316.65
$ws.Range("AE5").Value = 7565
$ws.Range("AF5").Value = 0.71
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 1.87
$ws.Range("AI5").Value = 551.27
$ws.Range("AJ5").Value = 43960757
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 2267
$ws.Range("E6").Value = 62
$ws.Range("F6").Value = 62
$ws.Range("G6").Value = 59
$ws.Range("H6").Value = -30
$ws.Range("I6").Value = -30
$ws.Range("K6").Value = 3601
$ws.Range("L6").Value = 564
$ws.Range("M6").Value = 3037
$ws.Range("N6").Value = 3037
$ws.Range("P6").Value = 220
$ws.Range("Q6").Value = 428
$ws.Range("R6").Value = -1128
$ws.Range("S6").Value = -44
$ws.Range("T6").Value = 266
$ws.Range("U6").Value = 162
$ws.Range("V6").Value = 133
$ws.Range("W6").Value = 2.72
$ws.Range("X6").Value = -1.31
$ws.Range("Y6").Value = -0.97
$ws.Range("Z6").Value = -0.82
$ws.Range("AA6").Value = 18.57
$ws.Range("AB6").Value = 1380.33
$ws.Range("AC6").Value = -68
$ws.Range("AD6").Value = -103.9
$ws.Range("AE6").Value = 7403
$ws.Range("AF6").Value = 0.95
$ws.Range("AI6").Value = -137.72
$ws.Range("AJ6").Value = 43960757
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Row 7
$ws.Range("D7").Value = 2550
$ws.Range("E7").Value = 136
$ws.Range("G7").Value = 151
$ws.Range("H7").Value = 111
$ws.Range("I7").Value = 112
$ws.Range("K7").Value = 3837
$ws.Range("L7").Value = 730
$ws.Range("M7").Value = 3107
$ws.Range("N7").Value = 3107
$ws.Range("P7").Value = 220
$ws.Range("Q7").Value = 209
$ws.Range("R7").Value = -150
$ws.Range("S7").Value = -42
$ws.Range("T7").Value = 159
$ws.Range("U7").Value = 23
$ws.Range("W7").Value = 5.34
$ws.Range("X7").Value = 4.36
$ws.Range("Y7").Value = 3.66
$ws.Range("Z7").Value = 2.99
$ws.Range("AA7").Value = 23.5
$ws.Range("AC7").Value = 256
$ws.Range("AD7").Value = 29.72
$ws.Range("AE7").Value = 7574
$ws.Range("AF7").Value = 1
$ws.Range("AG7").Value = 100
$ws.Range("AH7").Value = 1.32
$ws.Range("AI7").Value = 39.11

# Row 8
$ws.Range("D8").Value = 3454
$ws.Range("E8").Value = 278
$ws.Range("G8").Value = 295
$ws.Range("H8").Value = 222
$ws.Range("I8").Value = 222
$ws.Range("K8").Value = 4227
$ws.Range("L8").Value = 937
$ws.Range("M8").Value = 3291
$ws.Range("N8").Value = 3291
$ws.Range("P8").Value = 220
$ws.Range("Q8").Value = 266
$ws.Range("R8").Value = -247
$ws.Range("S8").Value = -38
$ws.Range("T8").Value = 112
$ws.Range("U8").Value = -52
$ws.Range("W8").Value = 8.050000000000001
$ws.Range("X8").Value = 6.42
$ws.Range("Y8").Value = 6.95
$ws.Range("Z8").Value = 5.5
$ws.Range("AA8").Value = 28.47
$ws.Range("AC8").Value = 506
$ws.Range("AD8").Value = 15.02
$ws.Range("AE8").Value = 8022
$ws.Range("AF8").Value = 0.95
$ws.Range("AG8").Value = 100
$ws.Range("AH8").Value = 1.32
$ws.Range("AI8").Value = 19.77

# Row 9
$ws.Range("D9").Value = 3711
$ws.Range("E9").Value = 335
$ws.Range("G9").Value = 351
$ws.Range("H9").Value = 268
$ws.Range("I9").Value = 268
$ws.Range("K9").Value = 4507
$ws.Range("L9").Value = 991
$ws.Range("M9").Value = 3516
$ws.Range("N9").Value = 3516
$ws.Range("P9").Value = 220
$ws.Range("Q9").Value = 282
$ws.Range("R9").Value = -235
$ws.Range("S9").Value = -39
$ws.Range("T9").Value = 116
$ws.Range("U9").Value = 120
$ws.Range("W9").Value = 9.029999999999999
$ws.Range("X9").Value = 7.22
$ws.Range("Y9").Value = 7.86
$ws.Range("Z9").Value = 6.14
$ws.Range("AA9").Value = 28.2
$ws.Range("AC9").Value = 609
$ws.Range("AD9").Value = 12.49
$ws.Range("AE9").Value = 8571
$ws.Range("AF9").Value = 0.89
$ws.Range("AG9").Value = 100
$ws.Range("AH9").Value = 1.32
$ws.Range("AI9").Value = 16.43
